$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.495.94"
$ws.Range("E2").Value = "  +1.20%  "
$ws.Range("D3").Value = "3.925.27"
$ws.Range("E3").Value = "  +1.23%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'484.83"
$ws.Range("E5").Value = "  +4.77%  "
$ws.Range("D6").Value = "'147.83"
$ws.Range("E6").Value = "  -0.42%  "
$ws.Range("D7").Value = "'0.621"
$ws.Range("E7").Value = "  -1.53%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "'0.724"
$ws.Range("E9").Value = "  -3.95%  "
$ws.Range("E10").Value = "  +6.96%  "
$ws.Range("D11").Value = "'0.0000357"
$ws.Range("E11").Value = "  +12.53%  "
$ws.Range("D12").Value = "'42.52"
$ws.Range("E12").Value = "  -3.83%  "
$ws.Range("D13").Value = "'10.52"
$ws.Range("E13").Value = "  +0.60%  "
$ws.Range("D14").Value = "4.553.47"
$ws.Range("E14").Value = "  +1.25%  "
$ws.Range("D15").Value = "'14.63"
$ws.Range("E15").Value = "  -0.94%  "
$ws.Range("D16").Value = "3.914.01"
$ws.Range("E16").Value = "  +0.84%  "
$ws.Range("E17").Value = "  -0.38%  "
$ws.Range("D18").Value = "'19.76"
$ws.Range("E18").Value = "  -1.86%  "
$ws.Range("E19").Value = "  -3.06%  "
$ws.Range("D20").Value = "68.621.40"
$ws.Range("E20").Value = "  +1.25%  "
$ws.Range("D21").Value = "'431.99"
$ws.Range("E21").Value = "  -0.06%  "
$ws.Range("B22").Value = "InternetComputer(DFINITY)"
$ws.Range("C22").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D22").Value = "'14.52"
$ws.Range("E22").Value = "  -2.70%  "
$ws.Range("B23").Value = "ImmutableX"
$ws.Range("C23").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D23").Value = "'3.34"
$ws.Range("E23").Value = "  +1.51%  "
$ws.Range("D24").Value = "'87.00"
$ws.Range("E24").Value = "  -1.75%  "
$ws.Range("D25").Value = "'11.37"
$ws.Range("E25").Value = "  +11.33%  "
$ws.Range("B26").Value = "RenderToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D26").Value = "'10.66"
$ws.Range("E26").Value = "  +3.13%  "
$ws.Range("B27").Value = "PancakeSwap"
$ws.Range("C27").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D27").Value = "'3.58"
$ws.Range("E27").Value = "  +0.72%  "
$ws.Range("D28").Value = "'38.16"
$ws.Range("E28").Value = "  +0.58%  "
$ws.Range("D29").Value = "'5.88"
$ws.Range("E29").Value = "  +6.73%  "
$ws.Range("D30").Value = "'709.83"
$ws.Range("E30").Value = "  -5.17%  "
$ws.Range("D31").Value = "'13.24"
$ws.Range("E31").Value = "  -4.01%  "
$ws.Range("E32").Value = "  -4.34%  "
$ws.Range("E33").Value = "  +2.71%  "
$ws.Range("D34").Value = "0.0₃0902"
$ws.Range("E34").Value = "  +31.85%  "
$ws.Range("D35").Value = "'41.49"
$ws.Range("E35").Value = "  -4.40%  "
$ws.Range("D36").Value = "'58.43"
$ws.Range("E36").Value = "  +1.50%  "
$ws.Range("D37").Value = "'0.151"
$ws.Range("E37").Value = "  -8.04%  "
$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").Value = "'5.52"
$ws.Range("E38").Value = "  -1.18%  "
$ws.Range("B39").Value = "Dai"
$ws.Range("C39").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D39").Value = "'0.999"
$ws.Range("E39").Value = "  -0.18%  "
$ws.Range("D40").Value = "'2.88"
$ws.Range("E40").Value = "  +9.78%  "
$ws.Range("D41").Value = "'0.0469"
$ws.Range("E41").Value = "  -2.41%  "
$ws.Range("D42").Value = "'3.04"
$ws.Range("E42").Value = "  +11.07%  "
$ws.Range("E43").Value = "  +1.75%  "
$ws.Range("D44").Value = "'0.346"
$ws.Range("E44").Value = "  -1.51%  "
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").Value = "'0.140"
$ws.Range("E45").Value = "  -0.90%  "
$ws.Range("B46").Value = "FirstDigitalUSD"
$ws.Range("C46").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D46").Value = "'1.00"
$ws.Range("E46").Value = "  +0.02%  "
$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D47").Value = "'2.17"
$ws.Range("E47").Value = "  +1.17%  "
$ws.Range("D48").Value = "'3.41"
$ws.Range("E48").Value = "  -1.24%  "
$ws.Range("B49").Value = "ApeXProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D49").Value = "'3.22"
$ws.Range("E49").Value = "  -2.12%  "
$ws.Range("D50").Value = "'147.69"
$ws.Range("E50").Value = "  +1.63%  "
$ws.Range("D51").Value = "'2.84"
$ws.Range("E51").Value = "  -2.48%  "
